$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq "Floor 0 & 1") {
        $cell.Value = "Floor 0-1"
    }
}

$ws.Range("C9").Select()
